# [version 27] UNet + GHCU Method: REGRESSION location and HEATMAP visibility
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 6 new rows before the old row 9 (try-1 ...), pushing the
# "try-*" block down from rows 9-13 to rows 15-19, keeping the blank
# separator row before it. ---
$ws.Range("8:13").Insert()

# --- Fill the newly inserted rows (8-13) with the new v21/v22-*/v23/v24
# REGRESSION location rows, matching the existing E/F/G(/H/I) layout. ---
$ws.Range("E8").Value = "v21"
$ws.Range("F8").Value = "epoch 100"
$ws.Range("G8").Value = 0.63815958611819401

$ws.Range("E9").Value = "v22"
$ws.Range("F9").Value = "epoch 100"
$ws.Range("G9").Value = 0.316436661235521

$ws.Range("E10").Value = "v22-2"
$ws.Range("F10").Value = "epoch 100"
$ws.Range("G10").Value = 0.31132784627410198
$ws.Range("H10").Value = 0.73397586585174701

$ws.Range("E11").Value = "v22-3"
$ws.Range("F11").Value = "epoch 100"
$ws.Range("G11").Value = 0.313779001121264

$ws.Range("E12").Value = "v23"
$ws.Range("F12").Value = "epoch 50"
$ws.Range("G12").Value = 0.34145764182721
$ws.Range("H12").Value = 0.73405422347594396
$ws.Range("I12").Value = "全猜visible"

$ws.Range("E13").Value = "v24"
$ws.Range("F13").Value = "epoch 50"
$ws.Range("G13").Value = 0.32692773974968298
$ws.Range("H13").Value = 0.51457451810061094

# --- Append two new try rows (try-6 / try-7) right after try-5 (now row 19). ---
$ws.Range("E20").Value = "try-6"
$ws.Range("F20").Value = "epoch 130"
$ws.Range("G20").Value = 0.32445426060138099

$ws.Range("E21").Value = "try-7"
$ws.Range("F21").Value = "epoch 140"
$ws.Range("G21").Value = 0.330743415082632

# --- Widen the new heatmap-visibility column (I). ---
$ws.Columns.Item(9).ColumnWidth = 12

# --- Restore the view/selection like the author left it. ---
$ws.Range("E9").Select()
